# Correção das notas do fórum para matc65 em 2021.2
# Zera as colunas B:J (visualizações diárias, total_views e nota_view)
# para os alunos que tinham pelo menos uma visualização registrada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Blocos contíguos de linhas que precisam ser zeradas (colunas B..J)
$ranges = @(
    "B3:J5",
    "B7:J8",
    "B10:J13",
    "B15:J19",
    "B21:J28",
    "B30:J30",
    "B32:J32",
    "B34:J36",
    "B38:J39",
    "B42:J44",
    "B47:J47",
    "B50:J50"
)

foreach ($rangeAddress in $ranges) {
    $ws.Range($rangeAddress).Value = 0
}
